# Regenerate save_data: update column G ("K") values for each game row.
# These values are externally computed (std/mean derived "s_vals"), so we
# simply overwrite the previously-stored numbers with the regenerated ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 2
    4  = 5
    5  = 2
    6  = 1
    7  = 0
    8  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 0
    26 = 1
    27 = 3
    28 = 0
    29 = 1
    30 = 2
    31 = 3
    32 = 0
    33 = 0
    34 = 1
    35 = 2
    36 = 2
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 2
    43 = 0
    44 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
